$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 4303.923
$ws.Range("J17").Value = 4478.6484
$ws.Range("L17").Value = 13435.9452
$ws.Range("N17").Value = -13771.9452
$ws.Range("H32").Value = 4190.1025
$ws.Range("I32").Value = 2000
$ws.Range("J32").Value = 4372.6113
$ws.Range("K32").Value = 2000
$ws.Range("L32").Value = 4372.6113
$ws.Range("M32").Value = -1674
$ws.Range("N32").Value = -5024.6113
$ws.Range("H62").Value = 9624.875
$ws.Range("J62").Value = 10199.866
$ws.Range("L62").Value = 10199.866
$ws.Range("N62").Value = -11447.866
$ws.Range("H64").Value = 6159.6772
$ws.Range("I64").Value = 5145.4546
$ws.Range("K64").Value = 5145.4546
$ws.Range("M64").Value = -4897.4546
$ws.Range("H65").Value = 9624.875
$ws.Range("J65").Value = 10199.866
$ws.Range("L65").Value = 50999.33
$ws.Range("N65").Value = -57239.33
$ws.Range("H67").Value = 6159.6772
$ws.Range("I67").Value = 5145.4546
$ws.Range("K67").Value = 5145.4546
$ws.Range("M67").Value = -4287.4546
$ws.Range("H74").Value = 6143.3335
$ws.Range("I74").Value = 3391.111
$ws.Range("K74").Value = 3391.111
$ws.Range("M74").Value = -2455.111
$ws.Range("H77").Value = 6143.3335
$ws.Range("I77").Value = 3391.111
$ws.Range("K77").Value = 16955.555
$ws.Range("M77").Value = -12275.555
$ws.Range("H98").Value = 1296
$ws.Range("I98").Value = 895.619
$ws.Range("J98").Value = 5500
$ws.Range("K98").Value = 895.619
$ws.Range("L98").Value = 5500
$ws.Range("M98").Value = 602.381
$ws.Range("N98").Value = -8496
$ws.Range("H107").Value = 2308.1
$ws.Range("I107").Value = 2308.1
$ws.Range("K107").Value = 2308.1
$ws.Range("M107").Value = -388.0999999999999
$ws.Range("H116").Value = 5436.6665
$ws.Range("I116").Value = 3965
$ws.Range("J116").Value = 6172.5
$ws.Range("K116").Value = 3965
$ws.Range("L116").Value = 6172.5
$ws.Range("M116").Value = -523
$ws.Range("N116").Value = -13056.5
$ws.Range("H122").Value = 1296
$ws.Range("I122").Value = 895.619
$ws.Range("J122").Value = 5500
$ws.Range("K122").Value = 2686.857
$ws.Range("L122").Value = 16500
$ws.Range("M122").Value = -236.857
$ws.Range("N122").Value = -21400
$ws.Range("H136").Value = 199769.67
$ws.Range("J136").Value = 199769.67
$ws.Range("L136").Value = 199769.67
$ws.Range("N136").Value = -209969.67
$ws.Range("H141").Value = 1545.9615
$ws.Range("I141").Value = 1256.6364
$ws.Range("J141").Value = 3137.25
$ws.Range("K141").Value = 3769.9092
$ws.Range("L141").Value = 9411.75
$ws.Range("M141").Value = 1410.0908
$ws.Range("N141").Value = -19771.75
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 50644.55
$ws.Range("I2").Value = 56119.5
$ws.Range("J2").Value = 1370
$ws.Range("K2").Value = 56119.5
$ws.Range("L2").Value = 1370
$ws.Range("M2").Value = -56006.5
$ws.Range("N2").Value = -1596
$ws.Range("H5").Value = 250
$ws.Range("I5").Value = 250
$ws.Range("K5").Value = 250
$ws.Range("M5").Value = -138
$ws.Range("H32").Value = 4324.6055
$ws.Range("I32").Value = 2828.3035
$ws.Range("K32").Value = 2828.3035
$ws.Range("M32").Value = -2541.3035
$ws.Range("H88").Value = 1266.6666
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 1266.6666
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 1266.6666
$ws.Range("M88").Value = $null
$ws.Range("N88").Value = -2078.6666
$ws.Range("H91").Value = 1266.6666
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 1266.6666
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 1266.6666
$ws.Range("M91").Value = $null
$ws.Range("N91").Value = -4074.6666
$ws.Range("H116").Value = 50644.55
$ws.Range("I116").Value = 56119.5
$ws.Range("J116").Value = 1370
$ws.Range("K116").Value = 56119.5
$ws.Range("L116").Value = 1370
$ws.Range("M116").Value = -53825.5
$ws.Range("N116").Value = -5958
$ws.Range("H138").Value = 57804.91
$ws.Range("J138").Value = 57804.91
$ws.Range("L138").Value = 57804.91
$ws.Range("N138").Value = -68084.91
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 50644.55
$ws.Range("I3").Value = 56119.5
$ws.Range("J3").Value = 1370
$ws.Range("K3").Value = 56119.5
$ws.Range("L3").Value = 1370
$ws.Range("M3").Value = -56005.5
$ws.Range("N3").Value = -1598
$ws.Range("H4").Value = 250
$ws.Range("I4").Value = 250
$ws.Range("K4").Value = 250
$ws.Range("M4").Value = -135
$ws.Range("H26").Value = 19461.8
$ws.Range("I26").Value = 10577.25
$ws.Range("J26").Value = 55000
$ws.Range("K26").Value = 10577.25
$ws.Range("L26").Value = 55000
$ws.Range("M26").Value = -10285.25
$ws.Range("N26").Value = -55584
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 4667.1904
$ws.Range("I58").Value = 4713.353
$ws.Range("J58").Value = 4471
$ws.Range("K58").Value = 4713.353
$ws.Range("L58").Value = 4471
$ws.Range("M58").Value = -4510.353
$ws.Range("N58").Value = -4877
$ws.Range("H132").Value = 1950.3
$ws.Range("I132").Value = 1950.3
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 5850.9
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -3320.9
$ws.Range("N132").Value = $null
$ws.Range("H134").Value = 35871.395
$ws.Range("I134").Value = 55830.117
$ws.Range("K134").Value = 167490.351
$ws.Range("M134").Value = -164955.351
$ws.Range("H136").Value = 4667.1904
$ws.Range("I136").Value = 4713.353
$ws.Range("J136").Value = 4471
$ws.Range("K136").Value = 14140.059
$ws.Range("L136").Value = 13413
$ws.Range("M136").Value = -11590.059
$ws.Range("N136").Value = -18513
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H28").Value = 3418.75
$ws.Range("I28").Value = 2925
$ws.Range("K28").Value = 8775
$ws.Range("M28").Value = -8543
$ws.Range("H56").Value = 20839172
$ws.Range("I56").Value = 20839172
$ws.Range("K56").Value = 20839172
$ws.Range("M56").Value = -20838642
$ws.Range("H97").Value = 222.27272
$ws.Range("I97").Value = 214.5
$ws.Range("J97").Value = 300
$ws.Range("K97").Value = 643.5
$ws.Range("L97").Value = 900
$ws.Range("M97").Value = -147.5
$ws.Range("N97").Value = -1892
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1951.3684
$ws.Range("I113").Value = 1061.2307
$ws.Range("K113").Value = 1061.2307
$ws.Range("M113").Value = 1108.7693
$ws.Range("H132").Value = 2838.7292
$ws.Range("I132").Value = 2360.875
$ws.Range("K132").Value = 7082.625
$ws.Range("M132").Value = -4552.625
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4680.067
$ws.Range("I7").Value = 2011.625
$ws.Range("K7").Value = 2011.625
$ws.Range("M7").Value = -1899.625
$ws.Range("H55").Value = 4142.385
$ws.Range("I55").Value = 3495.2
$ws.Range("J55").Value = 6299.6665
$ws.Range("K55").Value = 3495.2
$ws.Range("L55").Value = 6299.6665
$ws.Range("M55").Value = -3322.2
$ws.Range("N55").Value = -6645.6665
$ws.Range("H126").Value = 4680.067
$ws.Range("I126").Value = 2011.625
$ws.Range("K126").Value = 6034.875
$ws.Range("M126").Value = -3564.875
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 81794.62
$ws.Range("J41").Value = 81794.62
$ws.Range("L41").Value = 81794.62
$ws.Range("N41").Value = -82574.62
$ws.Range("H122").Value = 2331.125
$ws.Range("I122").Value = 1950.2858
$ws.Range("K122").Value = 5850.857400000001
$ws.Range("M122").Value = -3400.857400000001
